# ----------------------------------------------------------------------
# "fixed bugs related to capacity data"
#
# 1. Insert a brand-new worksheet "scenarios_new" as the first sheet,
#    containing a 27-row combinatorial scenario table (fast/base/slow
#    for Established / Battery / Hydrogen) with a 1/27 probability
#    column formatted as 0.000.
# 2. Move "variability" to be the second sheet (right after the new
#    scenarios_new sheet).
# 3. Tidy up leftover selections on a couple of the other sheets.
# ----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "scenarios_new" worksheet as the very first tab
# ------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($firstSheet)
$ws.Name = "scenarios_new"

# Column widths (best-fit widths captured from the authored workbook)
$ws.Columns.Item(3).ColumnWidth = 10.7109375
$ws.Columns.Item(4).ColumnWidth = 11.140625
$ws.Columns.Item(5).ColumnWidth = 7.42578125
$ws.Columns.Item(6).ColumnWidth = 9.5703125
$ws.Columns.Item(7).ColumnWidth = 12.28515625
$ws.Columns.Item(8).ColumnWidth = 20.140625
$ws.Columns.Item(9).ColumnWidth = 16.28515625
$ws.Columns.Item(10).ColumnWidth = 18.5703125
$ws.Columns.Item(11).ColumnWidth = 16.28515625

# Header row: Scenario | Name | Probability | Established | Battery | Hydrogen
$header = New-Object 'object[,]' 1,6
$header[0,0] = "Scenario"
$header[0,1] = "Name"
$header[0,2] = "Probability"
$header[0,3] = "Established"
$header[0,4] = "Battery"
$header[0,5] = "Hydrogen"
$ws.Range("A1:F1").Value = $header
$ws.Range("A1:F1").Font.Bold = $true
$ws.Range("G1").Font.Bold = $true

# Data rows 2-28: 27 combinations of fast/base/slow across three factors
$scenarioNames = @("OOO", "OOB", "OOP", "OBO", "OBB", "OBP", "OPO", "OPB", "OPP", "BOO", "BOB", "BOP", "BBO", "BBB", "BBP", "BPO", "BPB", "BPP", "POO", "POB", "POP", "PBO", "PBB", "PBP", "PPO", "PPB", "PPP")
$estCol  = @("fast", "fast", "fast", "fast", "fast", "fast", "fast", "fast", "fast", "base", "base", "base", "base", "base", "base", "base", "base", "base", "slow", "slow", "slow", "slow", "slow", "slow", "slow", "slow", "slow")
$battCol = @("fast", "fast", "fast", "base", "base", "base", "slow", "slow", "slow", "fast", "fast", "fast", "base", "base", "base", "slow", "slow", "slow", "fast", "fast", "fast", "base", "base", "base", "slow", "slow", "slow")
$hydCol  = @("fast", "base", "slow", "fast", "base", "slow", "fast", "base", "slow", "fast", "base", "slow", "fast", "base", "slow", "fast", "base", "slow", "fast", "base", "slow", "fast", "base", "slow", "fast", "base", "slow")

for ($i = 0; $i -lt 27; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $i
    $ws.Cells.Item($r, 2).Value = $scenarioNames[$i]
    $ws.Cells.Item($r, 4).Value = $estCol[$i]
    $ws.Cells.Item($r, 5).Value = $battCol[$i]
    $ws.Cells.Item($r, 6).Value = $hydCol[$i]
}

# Probability column: C2 on its own, C3:C28 as one shared-formula block
$ws.Range("C2").Formula = "=1/27"
$ws.Range("C3:C28").Formula = "=1/27"
$ws.Range("C2:C28").NumberFormat = "0.000"

$ws.Range("L11").Select()

# ------------------------------------------------------------------
# 2. Move "variability" so it sits right after "scenarios_new"
# ------------------------------------------------------------------
$variability = $wb.Worksheets.Item("variability")
$scenariosBase = $wb.Worksheets.Item("scenarios_base")
$variability.Move($scenariosBase)

# Re-fetch sheet references by name after Move(), since the prior
# object handles can end up pointing at the wrong tab once the sheet
# order has changed underneath them.
$ws = $wb.Worksheets.Item("scenarios_new")
$variability = $wb.Worksheets.Item("variability")
$scenariosBase = $wb.Worksheets.Item("scenarios_base")

$variability.Activate()
$variability.Range("A2").Select()

# ------------------------------------------------------------------
# 3. Selection tidy-up on "scenarios_base"
# ------------------------------------------------------------------
$scenariosBase.Activate()
$scenariosBase.Range("A1:K28").Select()

# ------------------------------------------------------------------
# Leave focus on the new first sheet, matching the authored workbook
# ------------------------------------------------------------------
$ws.Activate()
